$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.121.86'
$ws.Range("E2").Value = '  -1.79%  '

# Row 3
$ws.Range("D3").Value = '2.512.41'
$ws.Range("E3").Value = '  -2.18%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.25%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.583'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.35%  '

# Row 8
$ws.Range("E8").Value = '  +0.17%  '

# Row 9
$ws.Range("E9").Value = '  -2.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.92%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0804'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.61%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.41%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.112'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.56%  '

# Row 14
$ws.Range("D14").Value = '2.897.41'
$ws.Range("E14").Value = '  -2.17%  '

# Row 15
$ws.Range("D15").Value = '2.506.13'
$ws.Range("E15").Value = '  -2.36%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.97'
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.857'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.22%  '

# Row 18
$ws.Range("D18").Value = '42.226.18'
$ws.Range("E18").Value = '  -1.63%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.31%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0966'
$ws.Range("E20").Value = '  -3.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.65%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.72%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.45%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.61%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.73%  '

# Row 27
$ws.Range("E27").Value = '  +0.14%  '

# Row 28
$ws.Range("E28").Value = '  +9.44%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.62%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.42%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.07%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.17%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0782'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.02%  '

# Row 35
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.60'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.84%  '

# Row 36
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.09%  '

# Row 37
$ws.Range("E37").Value = '  -5.15%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.114'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.43%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.20%  '

# Row 40
$ws.Range("E40").Value = '  -1.16%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.35'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.05%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.80'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.54%  '

# Row 43
$ws.Range("E43").Value = '  +0.12%  '

# Row 44
$ws.Range("E44").Value = '  -2.70%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.031.63'
$ws.Range("E45").Value = '  -1.87%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0297'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.11%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.13%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.99%  '

# Row 49
$ws.Range("D49").Value = '2.761.82'
$ws.Range("E49").Value = '  -2.11%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.00%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.188'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.80%  '
